$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.246.66'
$ws.Range("E2").Value = '  +1.08%  '

$ws.Range("D3").Value = '1.617.46'
$ws.Range("E3").Value = '  +1.85%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.97'
$ws.Range("E5").Value = '  +0.80%  '

$ws.Range("E6").Value = '  +0.05%  '

$ws.Range("E7").Value = '  +0.93%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.247'
$ws.Range("E8").Value = '  +0.79%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0614'
$ws.Range("E9").Value = '  +0.84%  '

$ws.Range("E10").Value = '  +4.80%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0816'
$ws.Range("E11").Value = '  +0.93%  '

$ws.Range("D12").Value = '1.842.61'
$ws.Range("E12").Value = '  +1.88%  '

$ws.Range("D13").Value = '1.619.30'
$ws.Range("E13").Value = '  +1.74%  '

$ws.Range("E15").Value = '  +1.31%  '

$ws.Range("D16").Value = '26.263.27'
$ws.Range("E16").Value = '  +1.22%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.21'
$ws.Range("E17").Value = '  +3.77%  '

$ws.Range("E18").Value = '  +0.89%  '

$ws.Range("E19").Value = '  +0.01%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '200.70'
$ws.Range("E20").Value = '  +1.13%  '

$ws.Range("E21").Value = '  +1.57%  '

$ws.Range("E22").Value = '  +1.60%  '

$ws.Range("E23").Value = '  +1.15%  '

$ws.Range("E24").Value = '  +3.20%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.22'
$ws.Range("E25").Value = '  +1.01%  '

$ws.Range("E26").Value = '  +0.10%  '

$ws.Range("E27").Value = '  -1.02%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.13'
$ws.Range("E28").Value = '  +0.93%  '

$ws.Range("E29").Value = '  +1.81%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0513'
$ws.Range("E30").Value = '  +8.39%  '

$ws.Range("E31").Value = '  +0.86%  '

$ws.Range("E32").Value = '  +1.86%  '

$ws.Range("E33").Value = '  -0.43%  '

$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.41'
$ws.Range("E34").Value = '  +1.60%  '

$ws.Range("B35").Value = 'LidoDAOToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.49'

$ws.Range("D36").Value = '1.173.71'
$ws.Range("E36").Value = '  +4.73%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0164'
$ws.Range("E37").Value = '  +1.85%  '

$ws.Range("E38").Value = '  +2.74%  '

$ws.Range("E39").Value = '  +0.05%  '

$ws.Range("E40").Value = '  +0.07%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.494'
$ws.Range("E41").Value = '  +1.54%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.790'
$ws.Range("E42").Value = '  +1.43%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.33'
$ws.Range("E43").Value = '  +4.85%  '

$ws.Range("D44").Value = '1.754.88'
$ws.Range("E44").Value = '  +2.00%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.49'
$ws.Range("E45").Value = '  +0.78%  '

$ws.Range("E46").Value = '  +14.59%  '

$ws.Range("E47").Value = '  +3.05%  '

$ws.Range("E48").Value = '  +0.95%  '

$ws.Range("E49").Value = '  +1.12%  '

$ws.Range("E50").Value = '  +0.64%  '

$ws.Range("E51").Value = '  -0.26%  '
